$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Correct mislabeled rows: set label=1 and model to the correct source
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "gpt3"

$ws.Range("B23").Value = 1
$ws.Range("C23").Value = "gpt3"

$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "gpt3"

$ws.Range("B31").Value = 1
$ws.Range("C31").Value = "gemini"

# Reset column A width back to the standard default (drop the custom width)
$ws.Columns.Item(1).ColumnWidth = 8.43

# Update the view state (scroll position, zoom, selection) to match the
# final saved window state
$ws.Range("D66").Select()
$excel.ActiveWindow.Zoom = 160
